$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.527.78"
$ws.Range("E2").Value = "  -0.58%  "

$ws.Range("D3").Value = "'1.832.39"
$ws.Range("E3").Value = "  -0.67%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'312.05"
$ws.Range("E5").Value = "  -0.33%  "

$ws.Range("E6").Value = "  -0.07%  "

$ws.Range("D7").Value = "'0.4273"
$ws.Range("E7").Value = "  -0.32%  "

$ws.Range("D8").Value = "'0.3664"
$ws.Range("E8").Value = "  +0.63%  "

$ws.Range("D9").Value = "'0.07261"
$ws.Range("E9").Value = "  -0.86%  "

$ws.Range("D10").Value = "'0.8625"
$ws.Range("E10").Value = "  -1.80%  "

$ws.Range("D11").Value = "'20.60"
$ws.Range("E11").Value = "  -0.53%  "

$ws.Range("D12").Value = "'1.884.48"
$ws.Range("E12").Value = "  +1.11%  "

$ws.Range("D13").Value = "'5.395"
$ws.Range("E13").Value = "  +0.84%  "

$ws.Range("D14").Value = "'6.509"
$ws.Range("E14").Value = "  -0.17%  "

$ws.Range("D15").Value = "'0.06936"
$ws.Range("E15").Value = "  -0.25%  "

$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  -0.20%  "

$ws.Range("D17").Value = "'80.60"
$ws.Range("E17").Value = "  +1.31%  "

$ws.Range("E18").Value = "  -1.10%  "

$ws.Range("E19").Value = "  -0.16%  "

$ws.Range("D20").Value = "'15.38"
$ws.Range("E20").Value = "  -0.03%  "

$ws.Range("D21").Value = "'27.466.73"
$ws.Range("E21").Value = "  -0.97%  "

$ws.Range("D22").Value = "'5.154"
$ws.Range("E22").Value = "  +3.41%  "

$ws.Range("D23").Value = "'10.81"
$ws.Range("E23").Value = "  +4.85%  "

$ws.Range("D24").Value = "'2.056.90"
$ws.Range("E24").Value = "  -2.40%  "

$ws.Range("D25").Value = "'1.990"
$ws.Range("E25").Value = "  -0.04%  "

$ws.Range("D26").Value = "'154.50"
$ws.Range("E26").Value = "  -0.84%  "

$ws.Range("D27").Value = "'18.84"
$ws.Range("E27").Value = "  +1.44%  "

$ws.Range("D28").Value = "'5.098"
$ws.Range("E28").Value = "  -2.41%  "

$ws.Range("D29").Value = "'114.30"
$ws.Range("E29").Value = "  -4.58%  "

$ws.Range("D30").Value = "'1.819"
$ws.Range("E30").Value = "  -3.08%  "

$ws.Range("D31").Value = "'0.08851"
$ws.Range("E31").Value = "  -0.41%  "

$ws.Range("D32").Value = "'2.986"
$ws.Range("E32").Value = "  +1.02%  "

$ws.Range("D33").Value = "'0.7432"
$ws.Range("E33").Value = "  -1.36%  "

$ws.Range("D34").Value = "'4.534"
$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("D35").Value = "'1.129"
$ws.Range("E35").Value = "  +0.91%  "

$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("E37").Value = "  -1.40%  "

$ws.Range("D38").Value = "'0.05302"
$ws.Range("E38").Value = "  -2.44%  "

$ws.Range("D39").Value = "'0.01931"
$ws.Range("E39").Value = "  -0.20%  "

$ws.Range("D40").Value = "'2.796"
$ws.Range("E40").Value = "  -1.16%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.5070"
$ws.Range("E41").Value = "  -0.10%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.1656"
$ws.Range("E42").Value = "  -0.64%  "

$ws.Range("D43").Value = "'6.466"
$ws.Range("E43").Value = "  -2.10%  "

$ws.Range("D44").Value = "'8.296"
$ws.Range("E44").Value = "  -1.00%  "

$ws.Range("D45").Value = "'10.36"
$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("D46").Value = "'0.06478"
$ws.Range("E46").Value = "  -0.96%  "

$ws.Range("D47").Value = "'105.25"
$ws.Range("E47").Value = "  -0.74%  "

$ws.Range("D48").Value = "'0.4669"
$ws.Range("E48").Value = "  +0.41%  "

$ws.Range("D49").Value = "'1.000"
$ws.Range("E49").Value = "  -0.12%  "

$ws.Range("D50").Value = "'1.611"
$ws.Range("E50").Value = "  -1.59%  "

$ws.Range("E51").Value = "  -1.41%  "
